$d = $word.ActiveDocument

# The new bullet to add, right after the existing "singleton" bullet, inside
# the same "Design Patterns:" list (same ListParagraph style / numId).
$anchorText = "Made the Game class a singleton one."
$newText    = "Made balloon creation with Factory method."

# Find the paragraph whose text matches the anchor bullet. Paragraph.Range.Text
# includes the trailing paragraph mark (chr 13), so trim it before comparing.
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $ptext = $p.Range.Text.TrimEnd([char]13)
    if ($ptext -eq $anchorText) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate anchor paragraph '$anchorText'"
}

$target = $d.Paragraphs.Item($targetIndex)
$targetRange = $target.Range

# Paragraph.Range.End sits just past the paragraph mark, so back up one
# position to land right after the last visible character ("...one.") and
# before that mark (and before any bookmark, e.g. _GoBack, anchored there).
$insertPos = $targetRange.End - 1

# Insert the new bullet's text first, appending it to the existing run so it
# inherits that run's character formatting (the en-US language mark).
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertBefore($newText)

# Now split the paragraph: insert a paragraph mark right before the text we
# just added. This turns it into its own list paragraph (inheriting the same
# pPr - ListParagraph style / numId 2 list numbering) and leaves anything
# that followed the original insertion point (the _GoBack bookmark) attached
# to the new, final paragraph.
$splitPoint = $d.Range($insertPos, $insertPos)
$splitPoint.InsertBefore([char]13)
